$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I1 ("I0") and J1 ("IF"), matching the formatting of H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2..19 for columns I and J
$data = @{
    2  = @(4, 5)
    3  = @(3, 6)
    4  = @(5, 5)
    5  = @(6, 8)
    6  = @(7, 7)
    7  = @(9, 9)
    8  = @(3, 5)
    9  = @(4, 5)
    10 = @(5, 7)
    11 = @(4, 6)
    12 = @(5, 6)
    13 = @(5, 8)
    14 = @(8, 9)
    15 = @(6, 6)
    16 = @(6, 7)
    17 = @(6, 7)
    18 = @(4, 4)
    19 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
